# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.459.85"
$ws.Range("E2").Value = "  -1.92%  "
$ws.Range("D3").Value = "3.547.80"
$ws.Range("E3").Value = "  -2.76%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "'582.61"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").Value = "'173.05"
$ws.Range("E6").Value = "  -1.76%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.539.11"
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("E12").Value = "  -2.93%  "
$ws.Range("D13").Value = "'47.67"
$ws.Range("E14").Value = "  -4.65%  "
$ws.Range("D15").Value = "4.117.32"
$ws.Range("E15").Value = "  -2.87%  "
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("D17").Value = "'631.33"
$ws.Range("E17").Value = "  -6.11%  "
$ws.Range("D18").Value = "3.549.16"
$ws.Range("E18").Value = "  -3.34%  "
$ws.Range("D19").Value = "69.476.53"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "'17.48"
$ws.Range("E22").Value = "  -1.81%  "
$ws.Range("D23").Value = "'0.897"
$ws.Range("E23").Value = "  -3.89%  "
$ws.Range("D24").Value = "'16.09"
$ws.Range("E24").Value = "  -6.33%  "
$ws.Range("D25").Value = "'98.00"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("E28").Value = "  -5.29%  "
$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  -6.84%  "
$ws.Range("D30").Value = "'33.01"
$ws.Range("E30").Value = "  -5.37%  "
$ws.Range("E31").Value = "  -5.88%  "
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("E33").Value = "  -5.27%  "
$ws.Range("D34").Value = "'7.05"
$ws.Range("D35").Value = "'645.30"
$ws.Range("E35").Value = "  +11.20%  "
$ws.Range("E36").Value = "  -2.37%  "
$ws.Range("D37").Value = "'3.54"
$ws.Range("E37").Value = "  -11.24%  "
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("D39").Value = "'57.47"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "'0.0460"
$ws.Range("E41").Value = "  +0.65%  "
$ws.Range("E42").Value = "  -4.10%  "
$ws.Range("D43").Value = "3.408.68"
$ws.Range("E43").Value = "  -5.24%  "
$ws.Range("D44").Value = "'0.332"
$ws.Range("E44").Value = "  -4.63%  "
$ws.Range("E45").Value = "  -5.93%  "
$ws.Range("D46").Value = "'32.80"
$ws.Range("E46").Value = "  -6.11%  "
$ws.Range("E47").Value = "  -5.31%  "
$ws.Range("D48").Value = "'2.77"
$ws.Range("E48").Value = "  -3.80%  "
$ws.Range("E49").Value = "  -2.08%  "
$ws.Range("D50").Value = "'133.23"
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("E51").Value = "  +12.87%  "
